$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the sheet view back so A1 is the top-left visible cell again
# (instead of the previously scrolled-down A200), while keeping the
# existing selection on B218.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Give columns B:D (highest_infection_count, population,
# population_infected_percentage) explicit best-fit widths now that the
# queries/dashboard work is finished.
$ws.Columns.Item(2).ColumnWidth = 22.17
$ws.Columns.Item(3).ColumnWidth = 10.17
$ws.Columns.Item(4).ColumnWidth = 30

$ws.Range("B218").Select()
